# Weekly refresh of crypto price/volume figures (GitHub Actions bot).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.183.85"
$ws.Range("E2").Value = "'  +0.02%  "
$ws.Range("D3").Value = "'1.630.44"
$ws.Range("E3").Value = "'  -0.94%  "
$ws.Range("D5").Value = "'216.26"
$ws.Range("E5").Value = "'  -0.47%  "
$ws.Range("E8").Value = "'  -0.33%  "
$ws.Range("D9").Value = "'0.0624"
$ws.Range("E9").Value = "'  -0.85%  "
$ws.Range("D10").Value = "'20.34"
$ws.Range("E10").Value = "'  +2.03%  "
$ws.Range("D11").Value = "'0.0849"
$ws.Range("E11").Value = "'  +0.70%  "
$ws.Range("D12").Value = "'1.632.40"
$ws.Range("E12").Value = "'  -2.00%  "
$ws.Range("E13").Value = "'  +0.10%  "
$ws.Range("E14").Value = "'  +0.16%  "
$ws.Range("D15").Value = "'27.170.33"
$ws.Range("E15").Value = "'  +0.11%  "
$ws.Range("D16").Value = "'64.68"
$ws.Range("E16").Value = "'  -4.03%  "
$ws.Range("D17").Value = "'0.0₃0733"
$ws.Range("E17").Value = "'  -0.80%  "
$ws.Range("D18").Value = "'215.26"
$ws.Range("E18").Value = "'  -1.19%  "
$ws.Range("E19").Value = "'  +0.01%  "
$ws.Range("E20").Value = "'  +1.41%  "
$ws.Range("D21").Value = "'4.39"
$ws.Range("E21").Value = "'  -1.02%  "
$ws.Range("D22").Value = "'2.49"
$ws.Range("E22").Value = "'  -0.64%  "
$ws.Range("D23").Value = "'9.08"
$ws.Range("E23").Value = "'  -1.08%  "
$ws.Range("D24").Value = "'147.99"
$ws.Range("E24").Value = "'  +0.35%  "
$ws.Range("E25").Value = "'  +0.23%  "
$ws.Range("D26").Value = "'7.27"
$ws.Range("E26").Value = "'  -2.25%  "
$ws.Range("E27").Value = "'  -0.05%  "
$ws.Range("D28").Value = "'15.57"
$ws.Range("E28").Value = "'  -1.00%  "
$ws.Range("D29").Value = "'0.0504"
$ws.Range("E29").Value = "'  -0.05%  "
$ws.Range("E30").Value = "'  -0.47%  "
$ws.Range("E31").Value = "'  +0.31%  "
$ws.Range("E32").Value = "'  -0.62%  "
$ws.Range("D33").Value = "'1.313.82"
$ws.Range("E33").Value = "'  +3.75%  "
$ws.Range("B34").Value = "'HuobiToken"
$ws.Range("C34").Value = "'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").Value = "'2.53"
$ws.Range("E34").Value = "'  +3.24%  "
$ws.Range("B35").Value = "'LidoDAOToken"
$ws.Range("C35").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").Value = "'1.57"
$ws.Range("E35").Value = "'  -1.20%  "
$ws.Range("E36").Value = "'  -1.33%  "
$ws.Range("D37").Value = "'0.850"
$ws.Range("E37").Value = "'  +1.26%  "
$ws.Range("D38").Value = "'0.539"
$ws.Range("E38").Value = "'  -0.50%  "
$ws.Range("E39").Value = "'  +0.04%  "
$ws.Range("E40").Value = "'  +1.59%  "
$ws.Range("D41").Value = "'0.803"
$ws.Range("E41").Value = "'  -0.79%  "
$ws.Range("D42").Value = "'63.70"
$ws.Range("E42").Value = "'  +2.10%  "
$ws.Range("D43").Value = "'1.769.50"
$ws.Range("E43").Value = "'  -0.87%  "
$ws.Range("E44").Value = "'  -3.60%  "
$ws.Range("D45").Value = "'90.71"
$ws.Range("E45").Value = "'  -1.09%  "
$ws.Range("E46").Value = "'  -0.26%  "
$ws.Range("E47").Value = "'  -0.96%  "
$ws.Range("E48").Value = "'  +20.32%  "
$ws.Range("E49").Value = "'  +0.84%  "
$ws.Range("D50").Value = "'7.52"
$ws.Range("E50").Value = "'  -2.06%  "
$ws.Range("E51").Value = "'  -2.12%  "
